$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Column widths: widen column A, reset B:E toward the sheet's default width ---
$ws.Columns.Item(1).ColumnWidth = 20.833333333333332
$ws.Range("B1:E9").ColumnWidth = 7.833333333333333

# --- sheet view: zoom + selection moved to A10 ---
$ws.Application.ActiveWindow.Zoom = 100
$ws.Range("A10").Select()

# --- updated metric values ---
$ws.Range("C3").Value = 0.86580000000000001

$ws.Range("D8").Value = 0.22106191159593699
$ws.Range("E8").Value = 0.984708287874275

$ws.Range("C9").Value = 0.00000000000060524851590652497
$ws.Range("D9").Value = 0.49043076211320802
$ws.Range("E9").Value = 0.44895933125744703
